# Generate Report for Handback
#
# The localization-status report is regenerated periodically. For the
# "c771cdfd-9498-4d79-80a2-d7f330b13b6f" file, a new handback round-trip
# was produced, so the "Latest Handback DateTime" column (L) is refreshed
# with a newer timestamp on both locale sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("L3").Value = "2017-01-03 04:23:04"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("L3").Value = "2017-01-03 04:23:15"
